$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 118; existing rows 118..153 shift down to 119..154.
$ws.Rows("118:118").Insert()

# Populate the new row 118 with the weekly entry for Femacal de La Calera - Albahaca.
$ws.Range("A118").Value = 3
$ws.Range("B118").Value = "Femacal de La Calera"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = 44627
$ws.Range("E118").Value = 5
$ws.Range("F118").Value = 100112052
$ws.Range("G118").Value = "Albahaca"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 50
$ws.Range("K118").Value = 5000
$ws.Range("L118").Value = 5000
$ws.Range("M118").Value = 5000
$ws.Range("N118").Value = "$/docena de matas"
$ws.Range("O118").Value = "Provincia de Quillota"
$ws.Range("P118").Value = 833
$ws.Range("Q118").Value = 6
$ws.Range("R118").Value = "Hortaliza"
